$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove all existing hyperlinks on the sheet (B2, J2, J3, B3).
#    Calling Delete() on a Range-scoped Hyperlinks collection clears every
#    hyperlink on the worksheet, so one call is enough.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Update cell values for row 2 (new environment: ssurgwsoadev4-oci...)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B2").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"
$ws.Range("E2").Value = 11111003007
$ws.Range("G2").Value = "24/06/2022"
$ws.Range("J2").Value = "aseguradosgw@gmail.com"

# ---------------------------------------------------------------------------
# 3. Update cell values for row 3 (text content unchanged, just kept as-is)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Range("G3").Value = "29/10/2020"
$ws.Range("J3").Value = "aseguradosgw@gmail.com"

# ---------------------------------------------------------------------------
# 3b. Restore the quote-prefixed "General" number format on E2/G2/G3 that
#     setting .Value above resets back to the default style. F2 (untouched)
#     still carries that exact format, so copy it across via PasteSpecial.
# ---------------------------------------------------------------------------
$ws.Range("F2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Re-create the hyperlinks that remain in the edited workbook:
#      J2 -> mailto:aseguradosgw@gmail.com
#      J3 -> mailto:aseguradosgw@gmail.com
#      B3 -> https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do
#    Note: B2 intentionally no longer carries a hyperlink (value only).
# ---------------------------------------------------------------------------
$null = $ws.Hyperlinks.Add($ws.Range("J2"), "mailto:aseguradosgw@gmail.com", [Type]::Missing, [Type]::Missing, "aseguradosgw@gmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("J3"), "mailto:aseguradosgw@gmail.com", [Type]::Missing, [Type]::Missing, "aseguradosgw@gmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("B3"), "https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do", [Type]::Missing, [Type]::Missing, "https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")

# Re-assert the hyperlink cell style so it matches the original formatting.
$ws.Range("J2").Style = "Hipervínculo"
$ws.Range("J3").Style = "Hipervínculo"
$ws.Range("B3").Style = "Hipervínculo"

# ---------------------------------------------------------------------------
# 5. Update the active selection shown when the workbook is reopened.
# ---------------------------------------------------------------------------
$null = $ws.Range("B7").Select()
